$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2-27 is updated from serial date 45203
# (2023-10-04) to serial date 45204 (2023-10-05).
for ($row = 2; $row -le 27; $row++) {
    $ws.Cells.Item($row, 3).Value = 45204
}
